$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"

# Column B
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("B7").Value = "Wnt5a"

# Column C
$ws.Range("C2").Value = "Ror2"
$ws.Range("C3").Value = "Ror2"
$ws.Range("C4").Value = "Ror2"
$ws.Range("C5").Value = "Ror2"
$ws.Range("C6").Value = "Ror2"
$ws.Range("C7").Value = "Ror2"

# Column D
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "sCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "sCs"

# Column E
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 3

# Column F
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1

# Column G
$ws.Range("G2").Value = 0.1599003333333333
$ws.Range("G3").Value = 0.1599003333333333
$ws.Range("G4").Value = 0.1599003333333333
$ws.Range("G5").Value = 5.879152
$ws.Range("G6").Value = 5.879152
$ws.Range("G7").Value = 5.879152

# Column H
$ws.Range("H2").Value = 0.479701
$ws.Range("H3").Value = 0.479701
$ws.Range("H4").Value = 0.479701
$ws.Range("H5").Value = 17.637456
$ws.Range("H6").Value = 17.637456
$ws.Range("H7").Value = 17.637456

# Column I
$ws.Range("I2").Value = 0.0264777194346773
$ws.Range("I3").Value = 0.0264777194346773
$ws.Range("I4").Value = 0.0264777194346773
$ws.Range("I5").Value = 0.9735222805653226
$ws.Range("I6").Value = 0.9735222805653226
$ws.Range("I7").Value = 0.9735222805653226

# Column J
$ws.Range("J2").Value = 0.02647771943467731
$ws.Range("J3").Value = 0.02647771943467731
$ws.Range("J4").Value = 0.02647771943467731
$ws.Range("J5").Value = 0.9735222805653228
$ws.Range("J6").Value = 0.9735222805653228
$ws.Range("J7").Value = 0.9735222805653228

# Column K
$ws.Range("K2").Value = 2
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 2

# Column L
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 0.6666666666666666

# Column M
$ws.Range("M2").Value = 0.1378523333333333
$ws.Range("M3").Value = 7.692787333333334
$ws.Range("M4").Value = 0.1884283333333333
$ws.Range("M5").Value = 0.1378523333333333
$ws.Range("M6").Value = 7.692787333333334
$ws.Range("M7").Value = 0.1884283333333333

# Column N
$ws.Range("N2").Value = 0.413557
$ws.Range("N3").Value = 23.078362
$ws.Range("N4").Value = 0.565285
$ws.Range("N5").Value = 0.413557
$ws.Range("N6").Value = 23.078362
$ws.Range("N7").Value = 0.565285

# Column O
$ws.Range("O2").Value = 0.01719056794796269
$ws.Range("O3").Value = 0.9593118967607375
$ws.Range("O4").Value = 0.02349753529129985
$ws.Range("O5").Value = 0.01719056794796269
$ws.Range("O6").Value = 0.9593118967607375
$ws.Range("O7").Value = 0.02349753529129985

# Column P
$ws.Range("P2").Value = 0.01719056794796269
$ws.Range("P3").Value = 0.9593118967607375
$ws.Range("P4").Value = 0.02349753529129985
$ws.Range("P5").Value = 0.01719056794796269
$ws.Range("P6").Value = 0.9593118967607375
$ws.Range("P7").Value = 0.02349753529129985

# Column Q
$ws.Range("Q2").Value = 0.02204263405077778
$ws.Range("Q3").Value = 1.230079258862445
$ws.Range("Q4").Value = 0.03012975330944444
$ws.Range("Q5").Value = 0.8104548212213334
$ws.Range("Q6").Value = 45.22706603634134
$ws.Range("Q7").Value = 1.107798812773333

# Column R
$ws.Range("R2").Value = 0.198383706457
$ws.Range("R3").Value = 11.070713329762
$ws.Range("R4").Value = 0.271167779785
$ws.Range("R5").Value = 7.294093390992
$ws.Range("R6").Value = 407.043594327072
$ws.Range("R7").Value = 9.970189314960001

# Column S
$ws.Range("S2").Value = 0.0004551670350489126
$ws.Range("S3").Value = 0.02540039125277893
$ws.Range("S4").Value = 0.0006221611468494659
$ws.Range("S5").Value = 0.01673540091291378
$ws.Range("S6").Value = 0.9339115055079585
$ws.Range("S7").Value = 0.02287537414445039

# Column T
$ws.Range("T2").Value = 0.0004551670350489126
$ws.Range("T3").Value = 0.02540039125277893
$ws.Range("T4").Value = 0.0006221611468494661
$ws.Range("T5").Value = 0.01673540091291378
$ws.Range("T6").Value = 0.9339115055079587
$ws.Range("T7").Value = 0.02287537414445039
